$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 488 (pushes old 488:517 down to 490:519)
$ws.Rows("488:489").Insert()

# New row 488: Crimpson Seedless, Primera
$ws.Range("A488").Value2 = 5
$ws.Range("B488").Value2 = 'Macroferia Regional de Talca'
$ws.Range("C488").Value2 = 'Maule'
$ws.Range("D488").Value2 = 44783
$ws.Range("E488").Value2 = 7
$ws.Range("F488").Value2 = 'Fruta'
$ws.Range("G488").Value2 = 100109
$ws.Range("H488").Value2 = 'Uva'
$ws.Range("I488").Value2 = 100109001
$ws.Range("J488").Value2 = 'Uva'
$ws.Range("K488").Value2 = 'Crimpson Seedless'
$ws.Range("L488").Value2 = 'Primera'
$ws.Range("M488").Value2 = 250
$ws.Range("N488").Value2 = 9000
$ws.Range("O488").Value2 = 9000
$ws.Range("P488").Value2 = 9000
$ws.Range("Q488").Value2 = '$/bandeja 10 kilos'
$ws.Range("R488").Value2 = 'Provincia de Limarí'
$ws.Range("S488").Value2 = 900
$ws.Range("T488").Value2 = 10

# New row 489: Red Globe, Primera
$ws.Range("A489").Value2 = 5
$ws.Range("B489").Value2 = 'Macroferia Regional de Talca'
$ws.Range("C489").Value2 = 'Maule'
$ws.Range("D489").Value2 = 44783
$ws.Range("E489").Value2 = 7
$ws.Range("F489").Value2 = 'Fruta'
$ws.Range("G489").Value2 = 100109
$ws.Range("H489").Value2 = 'Uva'
$ws.Range("I489").Value2 = 100109001
$ws.Range("J489").Value2 = 'Uva'
$ws.Range("K489").Value2 = 'Red Globe'
$ws.Range("L489").Value2 = 'Primera'
$ws.Range("M489").Value2 = 150
$ws.Range("N489").Value2 = 9000
$ws.Range("O489").Value2 = 9000
$ws.Range("P489").Value2 = 9000
$ws.Range("Q489").Value2 = '$/bandeja 10 kilos'
$ws.Range("R489").Value2 = 'Provincia de Limarí'
$ws.Range("S489").Value2 = 900
$ws.Range("T489").Value2 = 10
